# Remove the trailing "20 min Erceg" note together with the blank
# paragraphs surrounding it, right after the last bullet point
# ("... durch die IDs realisiert" + the _GoBack bookmark), leaving the
# single blank paragraph that precedes the section break untouched.

$d = $word.ActiveDocument

# Locate the paragraph that carries the _GoBack bookmark (the last
# bullet of the "Logfile:" section) so the deletion is anchored to
# content rather than a hard-coded paragraph index.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*durch die IDs realisiert*") {
        $anchorIndex = $i
        break
    }
}

# The three paragraphs to remove are the ones immediately following the
# anchor: an empty paragraph, the "20 min ... Erceg" paragraph, and
# another empty paragraph.
$firstToDelete = $anchorIndex + 1
$lastToDelete = $anchorIndex + 3

$startPos = $d.Paragraphs($firstToDelete).Range.Start
$endPos = $d.Paragraphs($lastToDelete).Range.End

$r = $d.Range($startPos, $endPos)
$r.Delete()
